$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 6.176422
$ws.Range("H2").Value = 18.529266
$ws.Range("I2").Value = 0.3058063741187975
$ws.Range("J2").Value = 0.3058063741187975
$ws.Range("M2").Value = 62.19572466666667
$ws.Range("N2").Value = 186.587174
$ws.Range("O2").Value = 0.9009169178676326
$ws.Range("P2").Value = 0.9009169178676325
$ws.Range("Q2").Value = 384.1470421371426
$ws.Range("R2").Value = 3457.323379234284
$ws.Range("S2").Value = 0.2755061360353833
$ws.Range("T2").Value = 0.2755061360353832

$ws.Range("G3").Value = 6.176422
$ws.Range("H3").Value = 18.529266
$ws.Range("I3").Value = 0.3058063741187975
$ws.Range("J3").Value = 0.3058063741187975
$ws.Range("O3").Value = 0.03107099427955203
$ws.Range("P3").Value = 0.03107099427955203
$ws.Range("Q3").Value = 13.248536365596
$ws.Range("R3").Value = 119.236827290364
$ws.Range("S3").Value = 0.009501708100895706
$ws.Range("T3").Value = 0.009501708100895706

$ws.Range("G4").Value = 6.176422
$ws.Range("H4").Value = 18.529266
$ws.Range("I4").Value = 0.3058063741187975
$ws.Range("J4").Value = 0.3058063741187975
$ws.Range("O4").Value = 0.06801208785281536
$ws.Range("P4").Value = 0.06801208785281536
$ws.Range("Q4").Value = 29.00005745265533
$ws.Range("R4").Value = 261.000517073898
$ws.Range("S4").Value = 0.02079852998251858
$ws.Range("T4").Value = 0.02079852998251858

$ws.Range("I5").Value = 0.4631823009753332
$ws.Range("J5").Value = 0.4631823009753332
$ws.Range("M5").Value = 62.19572466666667
$ws.Range("N5").Value = 186.587174
$ws.Range("O5").Value = 0.9009169178676326
$ws.Range("P5").Value = 0.9009169178676325
$ws.Range("Q5").Value = 581.8391176530185
$ws.Range("R5").Value = 5236.552058877166
$ws.Range("S5").Value = 0.4172887710055354
$ws.Range("T5").Value = 0.4172887710055353

$ws.Range("I6").Value = 0.4631823009753332
$ws.Range("J6").Value = 0.4631823009753332
$ws.Range("O6").Value = 0.03107099427955203
$ws.Range("P6").Value = 0.03107099427955203
$ws.Range("S6").Value = 0.01439153462399433
$ws.Range("T6").Value = 0.01439153462399433

$ws.Range("I7").Value = 0.4631823009753332
$ws.Range("J7").Value = 0.4631823009753332
$ws.Range("O7").Value = 0.06801208785281536
$ws.Range("P7").Value = 0.06801208785281536
$ws.Range("S7").Value = 0.03150199534580353
$ws.Range("T7").Value = 0.03150199534580353

$ws.Range("I8").Value = 0.2310113249058692
$ws.Range("J8").Value = 0.2310113249058692
$ws.Range("M8").Value = 62.19572466666667
$ws.Range("N8").Value = 186.587174
$ws.Range("O8").Value = 0.9009169178676326
$ws.Range("P8").Value = 0.9009169178676325
$ws.Range("Q8").Value = 290.191195060892
$ws.Range("R8").Value = 2611.720755548028
$ws.Range("S8").Value = 0.208122010826714
$ws.Range("T8").Value = 0.2081220108267139

$ws.Range("I9").Value = 0.2310113249058692
$ws.Range("J9").Value = 0.2310113249058692
$ws.Range("O9").Value = 0.03107099427955203
$ws.Range("P9").Value = 0.03107099427955203
$ws.Range("R9").Value = 90.073522925388
$ws.Range("S9").Value = 0.007177751554661998
$ws.Range("T9").Value = 0.007177751554661998

$ws.Range("I10").Value = 0.2310113249058692
$ws.Range("J10").Value = 0.2310113249058692
$ws.Range("O10").Value = 0.06801208785281536
$ws.Range("P10").Value = 0.06801208785281536
$ws.Range("S10").Value = 0.01571156252449325
$ws.Range("T10").Value = 0.01571156252449325
